# Apply the edits described by the commit diff:
#  1. Shared string "Kinematic viscosity (cP)" -> "Dynamic viscosity (cP)"
#     (this is cell B1 on Sheet1, the column header)
#  2. Selected cell on Sheet1 changes from D10 to D8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the header text in B1 ("Kinematic viscosity (cP)" -> "Dynamic viscosity (cP)")
$ws.Range("B1").Value = "Dynamic viscosity (cP)"

# 2. Move/save the selection to D8 (was D10)
$ws.Range("D8").Select() | Out-Null
